$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.480.37'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '1.969.80'
$ws.Range('E3').Value = '  +3.58%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '326.49'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '0.4662'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = '0.3916'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = '46.22'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = '0.07930'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').Value = '0.9882'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('E12').Value = '  +4.40%  '
$ws.Range('D13').Value = '2.021.22'
$ws.Range('E13').Value = '  +6.20%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '5.802'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '0.07121'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = '87.77'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').Value = '0.000009931'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  +0.24%  '
$ws.Range('D22').Value = '29.518.88'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').Value = '5.536'
$ws.Range('E23').Value = '  +4.51%  '
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('D25').Value = '2.232.46'
$ws.Range('E25').Value = '  +4.65%  '
$ws.Range('D26').Value = '2.109'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '158.55'
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('D28').Value = '19.55'
$ws.Range('E28').Value = '  +0.77%  '
$ws.Range('D29').Value = '5.762'
$ws.Range('D30').Value = '119.66'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').Value = '1.901'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').Value = '0.09424'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').Value = '0.8914'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '5.237'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '1.321'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = '7.751'
$ws.Range('D41').Value = '0.5709'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '0.1794'
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('D43').Value = '9.649'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').Value = '0.000002908'
$ws.Range('E44').Value = '  +46.38%  '
$ws.Range('D45').Value = '2.735'
$ws.Range('E45').Value = '  +6.32%  '
$ws.Range('D46').Value = '11.78'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').Value = '0.5333'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('D48').Value = '2.159'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('D49').Value = '0.06922'
$ws.Range('E49').Value = '  -1.48%  '
$ws.Range('D50').Value = '1.831'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('D51').Value = '113.45'
$ws.Range('E51').Value = '  +0.23%  '
